# Generate Report for Handback
# Update status / error-detail cells to reflect a failed handback transform,
# and widen the "Error Detail" column so the new, longer messages are readable.

$wb = $excel.ActiveWorkbook

$statusMessage = "Handback transform failed"
$zhCnError = "Handback file name: 2nefoenn.gik is different with handoff file name: 0a1e432b-2776-4df2-84d5-90cb9e948489.46275d9a6ab67a20d1e05b214e361fbf0e8657d8.zh-cn."
$deDeError  = "Handback file name: 2nefoenn.gik is different with handoff file name: 0a1e432b-2776-4df2-84d5-90cb9e948489.46275d9a6ab67a20d1e05b214e361fbf0e8657d8.de-de."

# Excel's ColumnWidth (character units) is ~ raw OOXML column width minus 5/6,
# so to end up with a stored width of 40 we need to set ColumnWidth accordingly.
$targetColumnWidth = 40 - (5/6)

# --- Overview sheet: update the row for 0a1e432b-...-f99daa341c5b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusMessage
$wsOverview.Range("F3").Value = $statusMessage

# --- zh-cn sheet: update Status and Error Detail, widen column P ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusMessage
$wsZhCn.Range("P3").Value = $zhCnError
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColumnWidth

# --- de-de sheet: update Status and Error Detail, widen column P ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusMessage
$wsDeDe.Range("P3").Value = $deDeError
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColumnWidth
